# Add data for 2023-12-29
# Updates the 2023 (column J) running-total crime counts across the
# Citywide Totals, By Neighborhood, and per-neighborhood sheets to
# incorporate the day's newly reported incidents (some of which were
# backfilled into earlier report years, e.g. 2015/2020 columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7613
$ws.Range("J3").Value = 8006
$ws.Range("B4").Value = 1688
$ws.Range("G4").Value = 1477
$ws.Range("J4").Value = 1739
$ws.Range("J5").Value = 623
$ws.Range("J6").Value = 10942
$ws.Range("B7").Value = 23321
$ws.Range("G7").Value = 24701
$ws.Range("J7").Value = 28923

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 231
$ws.Range("J5").Value = 85
$ws.Range("J7").Value = 825
$ws.Range("J8").Value = 1830
$ws.Range("J11").Value = 527
$ws.Range("J13").Value = 34
$ws.Range("J15").Value = 355
$ws.Range("J18").Value = 233
$ws.Range("J19").Value = 843
$ws.Range("J20").Value = 627
$ws.Range("J21").Value = 84
$ws.Range("J23").Value = 264
$ws.Range("J25").Value = 152
$ws.Range("J29").Value = 1545
$ws.Range("J31").Value = 306
$ws.Range("J33").Value = 1304
$ws.Range("J34").Value = 132
$ws.Range("J37").Value = 889
$ws.Range("J42").Value = 1228
$ws.Range("J44").Value = 229
$ws.Range("J46").Value = 95
$ws.Range("J51").Value = 364
$ws.Range("J52").Value = 736
$ws.Range("J54").Value = 564
$ws.Range("B63").Value = 395
$ws.Range("G63").Value = 277
$ws.Range("J63").Value = 89
$ws.Range("J65").Value = 725
$ws.Range("J67").Value = 1049
$ws.Range("J68").Value = 61
$ws.Range("J73").Value = 287
$ws.Range("J75").Value = 86
$ws.Range("J76").Value = 409
$ws.Range("J77").Value = 202
$ws.Range("J79").Value = 795
$ws.Range("J83").Value = 587
$ws.Range("J85").Value = 1182
$ws.Range("J87").Value = 97
$ws.Range("J89").Value = 362
$ws.Range("J90").Value = 304
$ws.Range("J94").Value = 320
$ws.Range("J95").Value = 411
$ws.Range("J98").Value = 213
$ws.Range("J99").Value = 440
$ws.Range("B101").Value = 23321
$ws.Range("G101").Value = 24701
$ws.Range("J101").Value = 28923

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 250
$ws.Range("J6").Value = 262
$ws.Range("J7").Value = 825

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 144
$ws.Range("J7").Value = 527

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 109
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 362

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J6").Value = 338
$ws.Range("J7").Value = 1182

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 175
$ws.Range("J6").Value = 316
$ws.Range("J7").Value = 736

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J6").Value = 685
$ws.Range("J7").Value = 1830

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 587

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 294
$ws.Range("J3").Value = 433
$ws.Range("J7").Value = 1304

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 147
$ws.Range("J7").Value = 411

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 889

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 211
$ws.Range("J7").Value = 725

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 120
$ws.Range("J7").Value = 440

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 100
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 396
$ws.Range("J7").Value = 1049

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 143
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 564

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 468
$ws.Range("J3").Value = 544
$ws.Range("J5").Value = 57
$ws.Range("J6").Value = 393
$ws.Range("J7").Value = 1545

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 233
$ws.Range("J6").Value = 329
$ws.Range("J7").Value = 843

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 70
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 211
$ws.Range("J7").Value = 409

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J4").Value = 54
$ws.Range("J6").Value = 650
$ws.Range("J7").Value = 1228

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J5").Value = 16
$ws.Range("J6").Value = 34

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 73
$ws.Range("J7").Value = 264

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 230
$ws.Range("J3").Value = 267
$ws.Range("J7").Value = 795

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 208
$ws.Range("J7").Value = 627

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 61
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 74
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 55
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 108
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 364

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 66
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 97
